$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New StatQuery text (used for the "StatQuery" column C in rows 2-4).
# Replaces the previous apoc.text.split-based query with one that joins
# through an OPTIONAL MATCH on (p)<--(diag:diagnosis).
$newStatQuery = @"
MATCH (s:study)<--(p:participant)
OPTIONAL MATCH (p)<--(samp:sample)
MATCH (samp)<--(f:file)
WHERE f.experimental_strategy_and_data_subtypes in ['WXS']
OPTIONAL MATCH (p)<--(diag:diagnosis)
WITH DISTINCT samp,diag,s,p,f
RETURN
    count(distinct s) AS Studies,
    count(distinct p) AS Participants,
    count(distinct samp) AS Samples,
    count(distinct f) AS ``Files``
"@

# Update the StatQuery column (C) for the ParticipantsTab, SamplesTab and
# FilesTab rows so they all use the new query text.
$ws.Range("C2").Value2 = $newStatQuery
$ws.Range("C3").Value2 = $newStatQuery
$ws.Range("C4").Value2 = $newStatQuery

# The wrapped text got shorter, so the row heights shrink from 217 to 186.
$ws.Rows(2).RowHeight = 186
$ws.Rows(3).RowHeight = 186
$ws.Rows(4).RowHeight = 186

# Move the active selection down to B5 (scratch area below the table).
$ws.Range("B5").Select()
